$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "117/2022 Otsikko"
$ws.Range("B8").Value = "Valmis, Allekirjoitettu"
$ws.Range("C8").Value = "'1.1.2022"
$ws.Range("D8").Value = "Tweb"
$ws.Range("E8").Value = "Katu 1, 12345 Toimipaikka"
$ws.Range("F8").Value = "117/2022"
$ws.Range("G8").Value = "vastaanottaja Linna Yrjö"
$ws.Range("H8").Value = "Päätös"
$ws.Range("I8").Value = "123456789A"
$ws.Range("J8").Value = "Tyhjennysväli kielteinen"
$ws.Range("K8").Value = "'1.1.2021"
$ws.Range("L8").Value = "'3.10.2023"

$ws.Range("A3:P3").Copy()
$ws.Range("A8:P8").PasteSpecial(-4122)
